# The "VISUALIZATION" header in cell A1 of the TFEC sheet is renamed to
# "Type" (fixing a bug in the data-download header labels). The shared
# string "VISUALIZATION" is no longer referenced anywhere, and a new
# "Type" entry is used instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TFEC")

$ws.Range("A1").Value = "Type"

# Reset the active cell/selection back to A1 (top-left) instead of the
# stray A9 selection that was saved previously.
$ws.Range("A1").Select()
